$d = $word.ActiveDocument
$p1 = $d.Paragraphs.Item(1)
$r1 = $p1.Range
$r1.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item(2)
$r2 = $newPara.Range
$xmlSnippet = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>{</w:t></w:r><w:r><w:t>m</w:t></w:r><w:r><w:t>:</w:t></w:r><w:r><w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t>'dh</w:t></w:r><w:r><w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t>1</w:t></w:r><w:r><w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t>.gif'.asImage()</w:t></w:r><w:r><w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t>.fit(3</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t xml:space="preserve">00, </w:t></w:r><w:r><w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t>1</w:t></w:r><w:r><w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t>00)</w:t></w:r><w:r><w:t xml:space="preserve">}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r2.InsertXML($xmlSnippet)
$oldFieldPara = $d.Paragraphs.Item(3)
$oldFieldPara.Range.Delete()
